$wb = $excel.ActiveWorkbook

# --- Sheet: BASELINE_SIGLA ---
$ws1 = $wb.Worksheets.Item("BASELINE_SIGLA")

$ws1.Range("A2").Value = "MA0278"
$ws1.Range("A3").Value = "MA0279"
$ws1.Range("A4").Value = "MA0280"
$ws1.Range("B4").Value = "Homologation"
$ws1.Range("A5").Value = "MA0281"
$ws1.Range("B5").Value = "Pre-Production"
$ws1.Range("A6").Value = "MA0282"
$ws1.Range("B6").Value = "Homologation"

# Delete rows 7 and 8 (previously SCGTP0283CLD / SCGTP0284CLD)
$ws1.Rows("7:8").Delete()

# --- Sheet: Controle ---
$ws2 = $wb.Worksheets.Item("Controle")

$ws2.Range("A2").Value = "Homologation"
$ws2.Range("B2").Value = "4 x BACKUP 16GB RAM 8vCPU (WEBSERVER)"
$ws2.Range("G2").Value = "4"
$ws2.Range("H2").Value = "40"
$ws2.Range("I2").Value = "Hours/Week"
$ws2.Range("J2").Value = "On-Demand"
$ws2.Range("O2").Value = "2x Daily"
$ws2.Range("P2").Value = "10"

$ws2.Range("A3").Value = "Production"
$ws2.Range("B3").Value = "3 x BACKUP 32GB RAM 8vCPU (BACKUP)"
$ws2.Range("E3").Value = "m6id.2xlarge"
$ws2.Range("G3").Value = "3"
$ws2.Range("H3").Value = ""
$ws2.Range("I3").Value = "Always On"
$ws2.Range("J3").Value = "1 Yr No Upfront EC2 Instance Savings Plan"
$ws2.Range("L3").Value = "430"
$ws2.Range("O3").Value = "6x Daily"
$ws2.Range("P3").Value = "20"

$ws2.Range("B4").Value = "3 x BACKUP 32GB RAM 8vCPU (WEBSERVER)"
$ws2.Range("L4").Value = "800"

$ws2.Range("A5").Value = "Pre-Production"
$ws2.Range("B5").Value = "2 x BACKUP 16GB RAM 8vCPU (WEBSERVER)"
$ws2.Range("E5").Value = "c6i.2xlarge"
$ws2.Range("G5").Value = "2"
$ws2.Range("H5").Value = "40"
$ws2.Range("I5").Value = "Hours/Week"
$ws2.Range("J5").Value = "On-Demand"
$ws2.Range("L5").Value = "140"
$ws2.Range("O5").Value = "2x Daily"
$ws2.Range("P5").Value = "10"
